$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 05.02.2022 01:15"

# Update row 5 (Makro) with the new price data:
#  - B5 becomes the new current price
#  - C5 becomes the previous price (old B5 value)
#  - D5 becomes a text delta like "+0.4"
#  - E5 becomes a text timestamp instead of a numeric date serial
$ws.Range("B5").Value = 35.9
$ws.Range("C5").Value = 35.5

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "+0.4"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2022-02-05 01:16:37"
$ws.Range("E5").Style = "Normal"
